# Update the "想去人数" (interested-people count) values in the F column
# on both the "展览" and "全部类型" worksheets, which hold duplicate data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 479
    $ws.Range("F3").Value = 58
    $ws.Range("F4").Value = 30
}
